# Insert two new rows at 245-246 (pushes existing rows 245.. down to 247..)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A245:A246").EntireRow.Insert()

# New row 245: Ají, Inferno, Primera, Región de Arica y Parinacota
$ws.Range("A245").Value = 8
$ws.Range("B245").Value = "Terminal La Palmera de La Serena"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = 44855
$ws.Range("E245").Value = 4
$ws.Range("F245").Value = 100112021
$ws.Range("G245").Value = "Ají"
$ws.Range("H245").Value = "Inferno"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 480
$ws.Range("K245").Value = 18500
$ws.Range("L245").Value = 19000
$ws.Range("M245").Value = 18750
$ws.Range("N245").Value = "$/caja 10 kilos"
$ws.Range("O245").Value = "Región de Arica y Parinacota"
$ws.Range("P245").Value = 1875
$ws.Range("Q245").Value = 10
$ws.Range("R245").Value = "Hortaliza"

# New row 246: Ají, Inferno, Segunda, Región de Arica y Parinacota
$ws.Range("A246").Value = 8
$ws.Range("B246").Value = "Terminal La Palmera de La Serena"
$ws.Range("C246").Value = "Coquimbo"
$ws.Range("D246").Value = 44855
$ws.Range("E246").Value = 4
$ws.Range("F246").Value = 100112021
$ws.Range("G246").Value = "Ají"
$ws.Range("H246").Value = "Inferno"
$ws.Range("I246").Value = "Segunda"
$ws.Range("J246").Value = 400
$ws.Range("K246").Value = 12500
$ws.Range("L246").Value = 13000
$ws.Range("M246").Value = 12750
$ws.Range("N246").Value = "$/caja 10 kilos"
$ws.Range("O246").Value = "Región de Arica y Parinacota"
$ws.Range("P246").Value = 1275
$ws.Range("Q246").Value = 10
$ws.Range("R246").Value = "Hortaliza"
